$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# --- Registrar tiempos (horas consumidas) de las tareas completadas ---

# Fila 12: tarea ERS / Redacción Descripción general (Alonso) -> Hecho, 2h consumidas el Día 3 (col N)
$ws.Range("F12").Value = "Hecho"
$ws.Range("N12").Value = 2

# Fila 18: Administrar alumno / Diagrama de Robustez (Alonso) -> 2h consumidas el Día 5 (col T)
$ws.Range("T18").Value = 2

# Fila 19: Administrar profesor / Diagrama de Robustez (Alonso) -> 1h consumida el Día 5 (col T)
$ws.Range("T19").Value = 1

# Fila 20: Administrar grupo / Diagrama de Robustez (Alonso) -> 1h consumida el Día 5 (col T)
$ws.Range("T20").Value = 1

# Fila 27: ERS / Diagrama de Secuencia (Alonso) -> Hecho, 2h consumidas el Día 5 (col T)
$ws.Range("F27").Value = "Hecho"
$ws.Range("T27").Value = 2

# Recalcular la hoja para refrescar las formulas de horas restantes
$excel.CalculateFullRebuild()

# Actualizar la celda activa / seleccion visible en el panel inferior derecho
$ws.Range("N12").Select()
